$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header/label style (bold, bordered, centered) from B1
# onto the new column C header cell, then set its value.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C1").Value = "13-01-2023"

# Reordered fund rows (2-11), followed by avg/total (12-13). Column B keeps
# the values that were already present (just relocated); column C gets the
# new 13-01-2023 figures.
$data = @(
    @("Alpha Acciones",          9291.27,  9303.52),
    @("Alpha Mega",               1621.65,  1503.25),
    @("Bull Market",              5914.96,  5920.19),
    @("Compass Crecimiento",      0,        0),
    @("Delta Select",             0,        0),
    @("Fima Acciones",            0,        0),
    @("Fima PB Acciones",         0,        0),
    @("HF Acciones Argentinas",   158.17,   174.3),
    @("HF Acciones Lideres",      930.99,   931.88),
    @("Supergestion",             21680.31, 21273.02),
    @("avg",                      3959.73,  3910.62),
    @("total",                    39597.35, 39106.16)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value2 = $entry[1]
    $ws.Cells.Item($row, 3).Value2 = $entry[2]
    $row++
}
